$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
# ---------------------------------------------------------------------------
$titlePar = $d.Paragraphs(1)
$titlePar.Range.InsertParagraphAfter()

$metaPar = $d.Paragraphs(2)
$metaPar.Style = "Normal"

$metaRange = $metaPar.Range
$metaRange.Text = "Meta description: Experience the wild American West with Big Buffalo, a high-volatility slot game featuring 6 reels, 4,096 ways to win, and a Free Games feature. Play for free."

# Bold just the "Meta description" label (leave the rest of the sentence plain).
$labelStart = $metaRange.Start
$labelEnd = $labelStart + ("Meta description").Length
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Bold = 1

# Leave a leading (empty) run before the label, matching this document's
# existing body-paragraph convention of starting with a blank run.
$leadRange = $d.Range($metaRange.Start, $metaRange.Start)
$leadRange.InsertBefore("")

# ---------------------------------------------------------------------------
# 2) Drop the duplicate bold "Play Big Buffalo Free ..." paragraph near the
#    bottom of the document (paragraph 1 is the real title and must stay).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $par = $d.Paragraphs($i)
    if ($i -ne 1 -and $par.Range.Text.Contains("Play Big Buffalo Free - Exciting High-Variance Slot Game")) {
        $par.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Swap the final italic paragraph's copy for the DALLE image prompt text.
#    Scope the Find to the last paragraph only, so the identical sentence
#    that now also lives inside the meta-description paragraph is untouched.
# ---------------------------------------------------------------------------
$lastPar = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPar.Range

$oldBlurb = "Experience the wild American West with Big Buffalo, a high-volatility slot game featuring 6 reels, 4,096 ways to win, and a Free Games feature. Play for free."
# Placeholder ("|") stands in for the apostrophe: Find/Replace's ReplaceWith
# text gets smart-quoted on save, but a direct Range.Text character write
# does not, so the apostrophe is patched in afterwards.
$newPromptPlaceholder = "Prompt: DALLE, please create a feature image for the Big Buffalo slot game. The image should be in a cartoon style and should feature a happy Maya warrior with glasses. The warrior should be standing in front of a group of buffalo, with the mountains of the American West in the background. The image should be bright and colorful, with the warrior and buffalo being the main focus of attention. The image should also have the game|s logo prominently displayed somewhere in the design."

$found = $lastRange.Find.Execute($oldBlurb, $true, $false, $false, $false, $false, $true, 1, $false, $newPromptPlaceholder, 2)

if ($found) {
    $finalPar = $d.Paragraphs($d.Paragraphs.Count)
    $finalRange = $finalPar.Range
    $placeholderOffset = $finalRange.Text.IndexOf("|")
    if ($placeholderOffset -ge 0) {
        $placeholderStart = $finalRange.Start + $placeholderOffset
        $placeholderRange = $d.Range($placeholderStart, $placeholderStart + 1)
        $placeholderRange.Text = "'"
    }
}
